# daily auto push: 2026-02-24 07:14 UTC
# Insert a new daily-snapshot row at row 864 (date 2026/02/24, weekday 火,
# hour 15, ranking 201) which pushes the existing rows 864..905 down to
# 865..906, and updates the used range accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 864; everything below (864..905)
# shifts down to (865..906) automatically, just like pressing
# Ctrl+"+" on a selected row in the Excel UI.
$ws.Rows.Item(864).Insert()

# Column A holds dates formatted as plain text (e.g. "2026/12/29"), not
# real Excel date serials, so force the leading apostrophe to keep Excel
# from auto-converting the literal into a date value, then reset the
# cell style back to Normal so no stray "Quote Prefix" formatting is left
# behind on the cell (keeps the new row's styling identical to its
# neighbours, which also carry no explicit style).
$ws.Range("A864").Value = "'2026/02/24"
$ws.Range("A864").Style = "Normal"

$ws.Range("B864").Value = "火"
$ws.Range("C864").Value = 15
$ws.Range("D864").Value = 201
